# Fix: Home defaults to Latest Plan (Max Date) instead of Today
#
# Updates the two planning rows to the latest ("max date") plan data:
#  - shift the plan date forward one day
#  - update ticket ids to their newer values
#  - update patente (license plate) and cliente (client) values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 46020
$ws.Range("B2").Value = 413345733

# Row 3
$ws.Range("A3").Value = 46020
$ws.Range("B3").Value = 413235733

# Update cliente / patente text values.
# Order matters so the rebuilt shared-string table lines up the same way
# the original edit did: cliente (K) then patente (J), row 2 then row 3.
$ws.Range("K2").Value = "ECOTRANS"
$ws.Range("K3").Value = "AGRETRANS"
$ws.Range("J3").Value = "BSBJ91"
$ws.Range("J2").Value = "HCCR31"

# Reset the scroll position back to the top-left of the sheet (was parked
# at column H) now that the view no longer needs to default there.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1

# The ticket id values widened (7 -> 9 digits), so column B's best-fit
# width grows to fit the new content.
$ws.Columns.Item(2).ColumnWidth = 8.8
